$d = $word.ActiveDocument

$d.Content.Find.Execute(" is represented by composite primary keys (", $true, $false, $false, $false, $false,
                         $true, 1, $false, " replaces composite primary keys (", 2)

$d.Content.Find.Execute("is represented by composite primary keys (", $true, $false, $false, $false, $false,
                         $true, 1, $false, "replaces composite primary keys (", 2)
